# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it appears
#    (Overview!E2:F2, zh-cn!C2, de-de!C2 all share this string).
# 2) Narrow the "Status" columns (Overview E:F, zh-cn C, de-de C) from
#    ~17.22 chars to ~13.41 chars of column width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1) Update the status text ---------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2) Narrow the status columns -------------------------------------------
# ColumnWidth snaps to the host's internal pixel grid, so 12.5 is the input
# that lands closest to the target stored width (~13.41 characters).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
